# Auto update: 2025-12-05 02:00:49
# Applies the 2025-12-05 data refresh to the hedging/insurance analysis sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Date column (A2:A5) -------------------------------------------------
# The column stores the as-of date as plain text ("2025-12-03" -> "2025-12-05").
# Pre-formatting each cell as Text before assigning keeps Excel from
# auto-converting the date-like literal into a date serial number; clearing
# the formatting afterwards restores the cell's original (unstyled) look.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-12-05"
$ws.Range("A2").ClearFormats()

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2025-12-05"
$ws.Range("A3").ClearFormats()

$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2025-12-05"
$ws.Range("A4").ClearFormats()

$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2025-12-05"
$ws.Range("A5").ClearFormats()

# --- Row 2: UnitedHealth Group Incorporated / UNH ------------------------
$ws.Range("B2").Value = "UnitedHealth Group Incorporated"
$ws.Range("C2").Value = "UNH"
$ws.Range("D2").Value = 334.26
$ws.Range("E2").Value = 51.2
$ws.Range("F2").Value = 1.38
$ws.Range("G2").Value = 60
$ws.Range("H2").Value = 63
$ws.Range("I2").Value = 66
$ws.Range("J2").Value = 53
$ws.Range("K2").Value = 60.1
$ws.Range("L2").Value = "Pattern"
$ws.Range("M2").Value = "📈 매수 관찰 구간입니다."
$ws.Range("N2").Value = 52.43913937059539
$ws.Range("O2").Value = "⚪ 중립 구간"

# --- Row 3: MetLife, Inc. / MET ------------------------------------------
$ws.Range("B3").Value = "MetLife, Inc."
$ws.Range("C3").Value = "MET"
$ws.Range("D3").Value = 77.95999999999999
$ws.Range("E3").Value = 42.6
$ws.Range("F3").Value = 2.14
$ws.Range("G3").Value = 60
$ws.Range("H3").Value = 50
$ws.Range("I3").Value = 56
$ws.Range("J3").Value = 33
$ws.Range("K3").Value = 56.1
$ws.Range("L3").Value = "Pattern"
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 52.43913937059539
$ws.Range("O3").Value = "⚪ 중립 구간"

# --- Row 4: American International Group, I / AIG ------------------------
$ws.Range("B4").Value = "American International Group, I"
$ws.Range("C4").Value = "AIG"
$ws.Range("D4").Value = 77.45
$ws.Range("E4").Value = 44
$ws.Range("F4").Value = 1.85
$ws.Range("G4").Value = 60
$ws.Range("H4").Value = 43
$ws.Range("I4").Value = 40
$ws.Range("J4").Value = 40
$ws.Range("K4").Value = 49.7
$ws.Range("L4").Value = "Pattern"
$ws.Range("M4").Value = "⛔ 관망하십시오."
$ws.Range("N4").Value = 52.43913937059539
$ws.Range("O4").Value = "⚪ 중립 구간"

# --- Row 5: Prudential Financial, Inc. / PRU ------------------------------
$ws.Range("B5").Value = "Prudential Financial, Inc."
$ws.Range("C5").Value = "PRU"
$ws.Range("D5").Value = 110.55
$ws.Range("E5").Value = 64.40000000000001
$ws.Range("F5").Value = 2.39
$ws.Range("G5").Value = 60
$ws.Range("H5").Value = 43
$ws.Range("I5").Value = 36
$ws.Range("J5").Value = 43
$ws.Range("K5").Value = 48.1
$ws.Range("L5").Value = "Pattern"
$ws.Range("M5").Value = "⛔ 관망하십시오."
$ws.Range("N5").Value = 52.43913937059539
$ws.Range("O5").Value = "⚪ 중립 구간"
